# Changes done for Kaman new UI - header & footer

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("TC01_Verify_HomePage")
$wsData = $wb.Worksheets.Item("Testdata")

# --- Testdata sheet: add new element-type test data rows ---
# Values are written in this column/row order so that new shared-string
# entries get created in the same sequence as the target workbook
# (EleType1, EleType2, JSElement).
$wsData.Range("A6").Value = "EleType1"
$wsData.Range("A7").Value = "EleType2"
$wsData.Range("B6").Value = "JSElement"
$wsData.Range("B7").Value = "JSElement"

# Give the new rows the same bordered look used by the existing data rows,
# and bring the existing boolean cells (B3:B5) in line with it too.
$wsData.Range("A6:B7").Borders.LineStyle = 1
$wsData.Range("B3:B5").Borders.LineStyle = 1

# Update the selection shown when the sheet is active
$wsData.Range("A2:B7").Select() | Out-Null

# --- TC01_Verify_HomePage sheet ---
# Widen column C so the longer values fit
$wsMain.Columns.Item(3).ColumnWidth = 26.14

# Update the selection shown when the sheet is active
$wsMain.Select() | Out-Null
$wsMain.Range("A3:XFD7").Select() | Out-Null
